{"js": "// Replace the 25 division problems in the worksheet table with a new set\n// of values (same \"a\u00f7b=\" pattern), cell-for-cell, in document order.\n// Addressing by (row, col) avoids any ambiguity from values that are\n// shared between different cells (e.g. \"86\u00f72=\" is both an old value in\n// one cell and a new value written into another cell).\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Row indices that actually hold the five division problems in this\n// table (the intervening rows are blank \"answer\" rows).\nconst dataRows = [0, 4, 8, 12, 16];\n\n// Old -> new text for each of the 25 cells, in left-to-right,\n// top-to-bottom (document) order.\nconst replacements = [\n  \"40\u00f79=\", \"86\u00f72=\",\n  \"36\u00f77=\", \"11\u00f74=\",\n  \"19\u00f75=\", \"86\u00f75=\",\n  \"21\u00f78=\", \"18\u00f79=\",\n  \"82\u00f72=\", \"85\u00f78=\",\n  \"74\u00f74=\", \"80\u00f73=\",\n  \"89\u00f79=\", \"49\u00f72=\",\n  \"12\u00f75=\", \"43\u00f75=\",\n  \"71\u00f73=\", \"13\u00f77=\",\n  \"12\u00f76=\", \"31\u00f73=\",\n  \"59\u00f79=\", \"33\u00f78=\",\n  \"61\u00f78=\", \"98\u00f76=\",\n  \"28\u00f78=\", \"85\u00f76=\",\n  \"82\u00f76=\", \"83\u00f76=\",\n  \"36\u00f76=\", \"59\u00f73=\",\n  \"13\u00f72=\", \"92\u00f75=\",\n  \"74\u00f79=\", \"29\u00f77=\",\n  \"24\u00f74=\", \"32\u00f79=\",\n  \"74\u00f78=\", \"38\u00f76=\",\n  \"86\u00f72=\", \"60\u00f74=\",\n  \"85\u00f78=\", \"66\u00f76=\",\n  \"85\u00f74=\", \"93\u00f75=\",\n  \"88\u00f72=\", \"16\u00f72=\",\n  \"85\u00f77=\", \"33\u00f72=\",\n  \"18\u00f77=\", \"80\u00f76=\",\n];\n\n// Collect the 25 cell proxies first and load their current text so we can\n// sanity-check each one addresses the expected old value before writing\n// (cheap defensive check; addressing is purely positional so text search\n// ambiguity \u2014 e.g. \"86\u00f72=\" is both an old value and a value written\n// elsewhere \u2014 never comes into play).\nconst cells = [];\nfor (const row of dataRows) {\n  for (let col = 0; col < 5; col++) {\n    const cell = table.getCell(row, col);\n    cell.load(\"value\");\n    cells.push(cell);\n  }\n}\nawait context.sync();\n\nfor (let i = 0; i < cells.length; i++) {\n  const oldText = replacements[i * 2];\n  const newText = replacements[i * 2 + 1];\n  if (cells[i].value !== oldText) {\n    throw new Error(\n      `Unexpected cell text at index ${i}: expected \"${oldText}\", found \"${cells[i].value}\"`\n    );\n  }\n  cells[i].value = newText;\n}\n\nawait context.sync();\n", "ps1": "# Replace the 25 division problems in the worksheet table with a new set\n# of values (same \"a\u00f7b=\" pattern), cell-for-cell, in document order.\n# Addressing cells by their fixed (row, col) position avoids any ambiguity\n# from values that are shared between different cells (e.g. \"86\u00f72=\" is\n# both an old value in one cell and a new value written into another\n# cell later in the table).\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# Old -> new text for each of the 25 populated cells, in left-to-right,\n# top-to-bottom (document) order.\n$replacements = @(\n  \"40\u00f79=\", \"86\u00f72=\",\n  \"36\u00f77=\", \"11\u00f74=\",\n  \"19\u00f75=\", \"86\u00f75=\",\n  \"21\u00f78=\", \"18\u00f79=\",\n  \"82\u00f72=\", \"85\u00f78=\",\n  \"74\u00f74=\", \"80\u00f73=\",\n  \"89\u00f79=\", \"49\u00f72=\",\n  \"12\u00f75=\", \"43\u00f75=\",\n  \"71\u00f73=\", \"13\u00f77=\",\n  \"12\u00f76=\", \"31\u00f73=\",\n  \"59\u00f79=\", \"33\u00f78=\",\n  \"61\u00f78=\", \"98\u00f76=\",\n  \"28\u00f78=\", \"85\u00f76=\",\n  \"82\u00f76=\", \"83\u00f76=\",\n  \"36\u00f76=\", \"59\u00f73=\",\n  \"13\u00f72=\", \"92\u00f75=\",\n  \"74\u00f79=\", \"29\u00f77=\",\n  \"24\u00f74=\", \"32\u00f79=\",\n  \"74\u00f78=\", \"38\u00f76=\",\n  \"86\u00f72=\", \"60\u00f74=\",\n  \"85\u00f78=\", \"66\u00f76=\",\n  \"85\u00f74=\", \"93\u00f75=\",\n  \"88\u00f72=\", \"16\u00f72=\",\n  \"85\u00f77=\", \"33\u00f72=\",\n  \"18\u00f77=\", \"80\u00f76=\"\n)\n\n$rowCount = $t.Rows.Count\n$colCount = $t.Columns.Count\n$i = 0\n\nfor ($r = 1; $r -le $rowCount; $r++) {\n  # Skip the blank \"answer\" rows interleaved between the problem rows --\n  # only rows whose first cell holds text are addressed.\n  $probe = $t.Cell($r, 1).Range.Text\n  $probe = $probe -replace \"[\\x07\\x0d]\", \"\"\n  if ($probe -eq \"\") {\n    continue\n  }\n  for ($c = 1; $c -le $colCount; $c++) {\n    # Cell.Range.Text includes the trailing cell-mark characters (CR, BEL)\n    # -- strip them before comparing against the expected old value.\n    $oldText = $replacements[$i * 2]\n    $newText = $replacements[$i * 2 + 1]\n    $i++\n    $current = $t.Cell($r, $c).Range.Text -replace \"[\\x07\\x0d]\", \"\"\n    if ($current -ne $oldText) {\n      throw \"Unexpected cell text at row $r col ${c}: expected '$oldText', found '$current'\"\n    }\n    $t.Cell($r, $c).Range.Text = $newText\n  }\n}\n"}
